$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 139, shifting existing rows 139:165 down to 140:166
$ws.Range("A139:R139").EntireRow.Insert()

# Fill in the new row 139 with the new data record
$ws.Range("A139").Value = 7
$ws.Range("B139").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C139").Value = "Ñuble"
$ws.Range("D139").Value = 45173
$ws.Range("E139").Value = 16
$ws.Range("F139").Value = 100112031
$ws.Range("G139").Value = "Poroto verde"
$ws.Range("H139").Value = "Magnum"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 80
$ws.Range("K139").Value = 28000
$ws.Range("L139").Value = 28000
$ws.Range("M139").Value = 28000
$ws.Range("N139").Value = "$/malla 25 kilos"
$ws.Range("O139").Value = "Perú"
$ws.Range("P139").Value = 1120
$ws.Range("Q139").Value = 25
$ws.Range("R139").Value = "Hortaliza"

# Copy the date cell style (s="2") from the row below (now 140) into the new row 139
$ws.Range("D140").Copy()
$ws.Range("D139").PasteSpecial(-4122) # xlPasteFormats
